{"js": "// Apply the review-copy refresh described in the commit \"Added many more features\".\n// Each entry is an exact, unique (or intentionally repeated) sentence/phrase that is\n// replaced with its updated wording, mirroring the supplied OOXML diff.\nconst replacements = [\n  {\n    find: \"Play Gold Volcano Free: Exciting Cluster Slot Game\",\n    replace: \"Play Gold Volcano Free and Experience the Exciting Volcanic Slot Game\"\n  },\n  {\n    find: \"Cluster pays for wins with random features\",\n    replace: \"Exciting volcanic theme with catchy visual effects\"\n  },\n  {\n    find: \"RTP of 96.20% with medium-high volatility\",\n    replace: \"Random features and bonus features add excitement and potential wins\"\n  },\n  {\n    find: \"Well-designed graphic and sound effects\",\n    replace: \"Engaging soundtrack and sound effects enhance the gameplay experience\"\n  },\n  {\n    find: \"Great winning potential of up to 10,000x the bet\",\n    replace: \"High winning potential with an RTP of 96.20% and up to 10,000x multiplier\"\n  },\n  {\n    find: \"May not appeal to players who prefer traditional payline slots\",\n    replace: \"Medium-high volatility may not appeal to players seeking more frequent wins\"\n  },\n  {\n    find: \"Read our review of Gold Volcano, a cluster pay slot with random features and up to 10,000x wins. Play Gold Volcano for free at top online casinos.\",\n    replace: \"Read our review of Gold Volcano and play for free to enjoy the thrilling gameplay and high winning potential.\"\n  }\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the review-copy refresh described in the commit \"Added many more features\".\n# Each pair is an exact sentence/phrase replaced with its updated wording, mirroring\n# the supplied OOXML diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"Play Gold Volcano Free: Exciting Cluster Slot Game\"; Replace = \"Play Gold Volcano Free and Experience the Exciting Volcanic Slot Game\" },\n    @{ Find = \"Cluster pays for wins with random features\"; Replace = \"Exciting volcanic theme with catchy visual effects\" },\n    @{ Find = \"RTP of 96.20% with medium-high volatility\"; Replace = \"Random features and bonus features add excitement and potential wins\" },\n    @{ Find = \"Well-designed graphic and sound effects\"; Replace = \"Engaging soundtrack and sound effects enhance the gameplay experience\" },\n    @{ Find = \"Great winning potential of up to 10,000x the bet\"; Replace = \"High winning potential with an RTP of 96.20% and up to 10,000x multiplier\" },\n    @{ Find = \"May not appeal to players who prefer traditional payline slots\"; Replace = \"Medium-high volatility may not appeal to players seeking more frequent wins\" },\n    @{ Find = \"Read our review of Gold Volcano, a cluster pay slot with random features and up to 10,000x wins. Play Gold Volcano for free at top online casinos.\"; Replace = \"Read our review of Gold Volcano and play for free to enjoy the thrilling gameplay and high winning potential.\" }\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
